{"js": "// Update the date line and the 25 division problems in the practice table.\n// Every value is replaced positionally (row/column), not by text search,\n// because several problems share the same source text (e.g. \"61\u00f79=\" appears\n// twice) but map to different replacements.\n\nconst body = context.document.body;\n\n// 1) Date heading paragraph (\"2024-04-07 Sunday\" -> \"2024-04-08 Monday\").\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].getRange().insertText(\"2024-04-08 Monday\", \"Replace\");\n\n// 2) The 5x5 grid of division problems lives in table rows 0, 4, 8, 12, 16\n//    (the rows in between are blank spacer rows).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"58\u00f72=\", \"59\u00f79=\", \"29\u00f78=\", \"88\u00f75=\", \"43\u00f73=\"],\n  [\"46\u00f78=\", \"33\u00f79=\", \"59\u00f72=\", \"58\u00f73=\", \"71\u00f78=\"],\n  [\"76\u00f76=\", \"39\u00f76=\", \"11\u00f76=\", \"69\u00f75=\", \"71\u00f75=\"],\n  [\"48\u00f74=\", \"25\u00f77=\", \"14\u00f72=\", \"45\u00f75=\", \"51\u00f74=\"],\n  [\"91\u00f75=\", \"80\u00f75=\", \"32\u00f79=\", \"42\u00f74=\", \"98\u00f73=\"],\n];\nconst contentRows = [0, 4, 8, 12, 16];\n\nfor (let r = 0; r < contentRows.length; r++) {\n  const rowIndex = contentRows[r];\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(rowIndex, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n    cell.body.paragraphs.items[0].getRange().insertText(newValues[r][c], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice table.\n# Every value is replaced positionally (row/column), not by text search,\n# because several problems share the same source text (e.g. \"61\u00f79=\" appears\n# twice) but map to different replacements.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph (\"2024-04-07 Sunday\" -> \"2024-04-08 Monday\").\n$d.Paragraphs.Item(1).Range.Text = \"2024-04-08 Monday\"\n\n# 2) The 5x5 grid of division problems lives in table rows 1, 5, 9, 13, 17\n#    (1-based; the rows in between are blank spacer rows).\n$table = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"58\u00f72=\", \"59\u00f79=\", \"29\u00f78=\", \"88\u00f75=\", \"43\u00f73=\"),\n    @(\"46\u00f78=\", \"33\u00f79=\", \"59\u00f72=\", \"58\u00f73=\", \"71\u00f78=\"),\n    @(\"76\u00f76=\", \"39\u00f76=\", \"11\u00f76=\", \"69\u00f75=\", \"71\u00f75=\"),\n    @(\"48\u00f74=\", \"25\u00f77=\", \"14\u00f72=\", \"45\u00f75=\", \"51\u00f74=\"),\n    @(\"91\u00f75=\", \"80\u00f75=\", \"32\u00f79=\", \"42\u00f74=\", \"98\u00f73=\")\n)\n$contentRows = @(1, 5, 9, 13, 17)\n\nfor ($r = 0; $r -lt $contentRows.Length; $r++) {\n    $rowIndex = $contentRows[$r]\n    for ($c = 0; $c -lt 5; $c++) {\n        $cell = $table.Cell($rowIndex, $c + 1)\n        $cell.Range.Text = $newValues[$r][$c]\n    }\n}\n"}
